$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 47: "Local File Upload using Aggregation Service"
$ws.Range("A47").Value = "Local File Upload using Aggregation Service"
$ws.Range("B47").Value = "http://10.199.253.187:8085/api/cards/file"
$ws.Range("C47").Value = "POST"

# Make the URL in B47 a clickable hyperlink, like all the other endpoint cells.
$ws.Hyperlinks.Add($ws.Range("B47"), "http://10.199.253.187:8085/api/cards/file")

# Match formatting of the existing, similarly-styled rows above (A45/B45/C45,
# D44) so the new row reuses the same cell styles rather than creating new
# duplicate style entries. Done after adding the hyperlink so the explicit
# formatting below wins over the auto-applied hyperlink style.
$ws.Range("A45").Copy()
$ws.Range("A47").PasteSpecial(-4122)

$ws.Range("B45").Copy()
$ws.Range("B47").PasteSpecial(-4122)

$ws.Range("C45").Copy()
$ws.Range("C47").PasteSpecial(-4122)

$ws.Range("D44").Copy()
$ws.Range("D47").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Match the row height used by the other data rows.
$ws.Rows.Item(47).RowHeight = 15

# The new row becomes the active selection, as in the authored workbook.
$ws.Range("A47").Select()
